$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C (old C shifts right to become E)
$ws.Range("C1:D1").EntireColumn.Insert()

# Set header row values - each week's column keeps its own date label,
# newest week now lives in B, and the older weeks shift right
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = "Jun_13"

# Fill the two new columns (C and D) with "UN" for data rows 2-27
$ws.Range("C2:D27").Value = "UN"

# Keep the same custom column width (8.0) applied to the newly inserted
# columns, matching the formatting of the rest of the table
$ws.Columns("C").ColumnWidth = 7.1666666
$ws.Columns("D").ColumnWidth = 7.1666666
$ws.Columns("E").ColumnWidth = 7.1666666
